$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 34: new array formula ATAN2(B34:C34,D34) spilled over T34:U34 ---
$ws.Range("T34:U34").FormulaArray = "=ATAN2(B34:C34,D34)"

# --- Row 48: new array formula SUM(MOD(B48:E48,2)) in T48, with a dedicated font ---
$ws.Range("T48").FormulaArray = "=SUM(MOD(B48:E48,2))"
$ws.Range("T48").Font.Size = 10
$ws.Range("T48").Font.Name = "Arial Unicode MS"
$ws.Rows.Item(48).RowHeight = 16.5

# --- Row 50: new SUMPRODUCT example row ---
$ws.Range("A50").Value = "SUMPRODUCT"
$ws.Range("B50").Value = 0
$ws.Range("C50").Value = 1
$ws.Range("D50").Value = 2
$ws.Range("E50").Value = 3
$ws.Range("F50").Value = "A"
$ws.Range("I50").Formula = "=SUMPRODUCT(B50:C50,C50:D50,D50:E50)"
$ws.Range("J50").Formula = "=SUMPRODUCT(C50:D50,D50:E50,E50:F50)"
$ws.Range("K50").Formula = "=SUMPRODUCT(D50:E50,E50:F50,F50:G50)"
$ws.Range("L50").Formula = "=SUMPRODUCT(E50:F50,F50:G50,G50)"

# --- Restore the selection as left by the editor ---
$ws.Range("K51").Select()
